$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "28.096.21"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -0.54%  "

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.876.27"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -1.90%  "

$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.46%  "

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "313.45"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.04%  "

$ws.Cells.Item(6, 5).Value = "  +0.38%  "

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.5098"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.41%  "

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.3845"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -2.52%  "

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.09098"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -2.81%  "

$ws.Cells.Item(10, 5).Value = "  -1.74%  "

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "41.57"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -0.77%  "

$ws.Cells.Item(12, 5).Value = "  -0.83%  "

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "20.75"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -0.84%  "

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.875.28"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -1.30%  "

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.209"
$cell.NumberFormat = "General"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.003"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +0.51%  "

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.00001113"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -1.21%  "

$ws.Cells.Item(18, 5).Value = "  -1.72%  "

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.06596"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.25%  "

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "18.15"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +1.22%  "

$ws.Cells.Item(21, 5).Value = "  +0.31%  "

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.108"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -1.90%  "

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "28.112.05"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.61%  "

$ws.Cells.Item(24, 5).Value = "  +0.12%  "

$ws.Cells.Item(25, 5).Value = "  -1.41%  "

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.092.57"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -1.25%  "

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.542"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -3.34%  "

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "20.81"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -1.18%  "

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "157.78"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.23%  "

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "126.65"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -0.55%  "

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.064"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -3.01%  "

$ws.Cells.Item(32, 5).Value = "  -1.76%  "

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.611"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.47%  "

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.600"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.24%  "

$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "9.691"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +0.08%  "

$ws.Cells.Item(36, 2).Value = "VeChain"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.02433"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +0.66%  "

$ws.Cells.Item(37, 2).Value = "Hedera"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.06567"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -2.01%  "

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.2178"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -0.47%  "

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.211"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -3.62%  "

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.261"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +0.74%  "

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.57"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.34%  "

$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.6401"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.12%  "

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.913"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -1.70%  "

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "13.24"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.97%  "

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.6017"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.10%  "

$ws.Cells.Item(46, 5).Value = "  -0.94%  "

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.276"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.01%  "

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.236"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +4.55%  "

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.998"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -1.20%  "

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "121.33"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -1.50%  "

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "79.76"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +1.46%  "
